$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = "'44.182.57"
$ws.Range('E2').Value = '  +2.37%  '

# Row 3
$ws.Range('D3').Value = "'2.277.52"
$ws.Range('E3').Value = '  +2.75%  '

# Row 4
$ws.Range('E4').Value = '  -0.42%  '

# Row 5
$ws.Range('D5').Value = "'318.37"
$ws.Range('E5').Value = '  +0.79%  '

# Row 6
$ws.Range('D6').Value = "'105.87"
$ws.Range('E6').Value = '  +7.64%  '

# Row 7
$ws.Range('E7').Value = '  +1.22%  '

# Row 8
$ws.Range('E8').Value = '  -0.35%  '

# Row 9
$ws.Range('D9').Value = "'0.573"
$ws.Range('E9').Value = '  +2.24%  '

# Row 10
$ws.Range('D10').Value = "'39.04"
$ws.Range('E10').Value = '  +7.01%  '

# Row 11
$ws.Range('E11').Value = '  +1.71%  '

# Row 12
$ws.Range('E12').Value = '  +1.91%  '

# Row 13
$ws.Range('E13').Value = '  +1.72%  '

# Row 14
$ws.Range('D14').Value = "'2.628.55"
$ws.Range('E14').Value = '  +2.87%  '

# Row 15
$ws.Range('D15').Value = "'0.882"
$ws.Range('E15').Value = '  +2.32%  '

# Row 16
$ws.Range('E16').Value = '  +3.38%  '

# Row 17
$ws.Range('D17').Value = "'2.287.13"
$ws.Range('E17').Value = '  +3.47%  '

# Row 18
$ws.Range('D18').Value = "'44.100.49"
$ws.Range('E18').Value = '  +2.79%  '

# Row 19
$ws.Range('D19').Value = "'14.15"
$ws.Range('E19').Value = '  -4.13%  '

# Row 20
$ws.Range('D20').Value = "'0.0000100"
$ws.Range('E20').Value = '  +4.08%  '

# Row 21
$ws.Range('D21').Value = "'6.56"
$ws.Range('E21').Value = '  +2.68%  '

# Row 22
$ws.Range('D22').Value = "'66.24"
$ws.Range('E22').Value = '  +1.55%  '

# Row 23
$ws.Range('E23').Value = '  +1.95%  '

# Row 24
$ws.Range('D24').Value = "'237.87"
$ws.Range('E24').Value = '  +0.56%  '

# Row 25
$ws.Range('E25').Value = '  +4.59%  '

# Row 26
$ws.Range('D26').Value = "'1.00"
$ws.Range('E26').Value = '  -0.10%  '

# Row 27
$ws.Range('D27').Value = "'10.28"
$ws.Range('E27').Value = '  +2.19%  '

# Row 28
$ws.Range('D28').Value = "'39.11"
$ws.Range('E28').Value = '  +15.52%  '

# Row 29
$ws.Range('E29').Value = '  -0.19%  '

# Row 30
$ws.Range('D30').Value = "'6.54"
$ws.Range('E30').Value = '  +4.09%  '

# Row 31
$ws.Range('D31').Value = "'163.93"
$ws.Range('E31').Value = '  +4.63%  '

# Row 32
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').Value = "'20.56"
$ws.Range('E32').Value = '  +0.47%  '

# Row 33
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = "'0.0884"
$ws.Range('E33').Value = '  +1.47%  '

# Row 34
$ws.Range('E34').Value = '  -0.94%  '

# Row 35
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').Value = "'3.27"
$ws.Range('E35').Value = '  +1.22%  '

# Row 36
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').Value = "'2.08"
$ws.Range('E36').Value = '  +3.72%  '

# Row 37
$ws.Range('E37').Value = '  +13.45%  '

# Row 38
$ws.Range('E38').Value = '  -0.61%  '

# Row 39
$ws.Range('B39').Value = 'NEARProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D39').Value = "'3.95"
$ws.Range('E39').Value = '  +7.62%  '

# Row 40
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').Value = "'4.50"
$ws.Range('E40').Value = '  +1.28%  '

# Row 41
$ws.Range('D41').Value = "'0.0326"
$ws.Range('E41').Value = '  +0.42%  '

# Row 42
$ws.Range('D42').Value = "'15.42"
$ws.Range('E42').Value = '  +27.20%  '

# Row 43
$ws.Range('E43').Value = '  -0.48%  '

# Row 44
$ws.Range('D44').Value = "'1.765.53"
$ws.Range('E44').Value = '  -6.92%  '

# Row 45
$ws.Range('E45').Value = '  +0.49%  '

# Row 46
$ws.Range('D46').Value = "'85.75"

# Row 47
$ws.Range('E47').Value = '  -0.44%  '

# Row 48
$ws.Range('B48').Value = 'ordi'
$ws.Range('C48').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D48').Value = "'75.54"
$ws.Range('E48').Value = '  +0.66%  '

# Row 49
$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').Value = "'8.83"
$ws.Range('E49').Value = '  +3.09%  '

# Row 50
$ws.Range('D50').Value = "'59.64"
$ws.Range('E50').Value = '  -1.44%  '

# Row 51
$ws.Range('D51').Value = "'104.48"
$ws.Range('E51').Value = '  +3.18%  '

